# mlResults/ML Model Metadata.xlsx - add "Model2 from 2020 to 2021" column + row cleanup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column E (old E:J -> new F:K)
$ws.Range("E1").EntireColumn.Insert()

# 2. New column header + numeric index values (mirrors columns B/C/D)
$ws.Range("E1").Value = "Unnamed: 0.1.1.1"
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 4

# 3. The old duplicate "Model 2 / crossectional random forest" row (row 7)
#    is dropped, and the second duplicate (row 8, with corrected wording)
#    slides up to take its place as row 7.
$ws.Range("A7:K7").EntireRow.Delete()
$ws.Range("A7").Value = 5
$ws.Range("I7").Value = "run random forest model year by year bbasis"
